# edit.ps1 - Applies the "New crime data collected" weekly update to the
# 76th Precinct CompStat workbook.
#
# 1) Updates the two rich-text header captions (issue number + reporting
#    week dates) by rewriting only the specific substring run, preserving
#    the rest of the caption text.
# 2) Updates the weekly crime-statistics grid (rows 16-30, columns C:N)
#    with the newly collected figures. A few cells flip from the "no data"
#    placeholder text (shared strings "0" / "***.*") to real numbers, so
#    their number format is set explicitly to match sibling cells in the
#    same column (integer count format or the 1-decimal % format).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header captions -------------------------------------------------

# A8: "Volume 31   Number  39" -> "...  40"
$ws.Range("A8").Characters(21, 2).Text = "40"

# C9: "Report Covering the Week  9/23/2024  Through  9/29/2024"
#  -> "...  9/30/2024  Through  10/6/2024"
$ws.Range("C9").Characters(27, 9).Text = "9/30/2024"
$ws.Range("C9").Characters(47, 9).Text = "10/6/2024"

# --- Weekly crime statistics grid ------------------------------------

$ws.Range("C16").Value = 1
$ws.Range("C16").NumberFormat = "#,##0"
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("G16").Value = 4
$ws.Range("H16").Value = -25
$ws.Range("I16").Value = 49
$ws.Range("J16").Value = 51
$ws.Range("K16").Value = -3.921568627450
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -37.179487179487
$ws.Range("N16").Value = -82.5
$ws.Range("F17").Value = 11
$ws.Range("G17").Value = 6
$ws.Range("H17").Value = 83.333333333333
$ws.Range("I17").Value = 95
$ws.Range("K17").Value = 13.095238095238
$ws.Range("L17").Value = 26.666666666666
$ws.Range("M17").Value = 25
$ws.Range("N17").Value = -60.084033613445
$ws.Range("F18").Value = 5
$ws.Range("H18").Value = 25
$ws.Range("I18").Value = 71
$ws.Range("K18").Value = -4.054054054054
$ws.Range("L18").Value = -15.476190476190
$ws.Range("M18").Value = -16.470588235294
$ws.Range("N18").Value = -78.353658536585
$ws.Range("C19").Value = 1
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = -75
$ws.Range("F19").Value = 16
$ws.Range("G19").Value = 14
$ws.Range("H19").Value = 14.285714285714
$ws.Range("I19").Value = 157
$ws.Range("J19").Value = 126
$ws.Range("K19").Value = 24.603174603174
$ws.Range("L19").Value = -0.632911392405
$ws.Range("M19").Value = -10.285714285714
$ws.Range("N19").Value = -13.259668508287
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 2
$ws.Range("D20").NumberFormat = "#,##0"
$ws.Range("E20").Value = 50
$ws.Range("E20").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F20").Value = 6
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 100
$ws.Range("I20").Value = 48
$ws.Range("J20").Value = 49
$ws.Range("K20").Value = -2.040816326530
$ws.Range("L20").Value = 9.090909090909
$ws.Range("M20").Value = -2.040816326530
$ws.Range("N20").Value = -84.516129032258
$ws.Range("C21").Value = 8
$ws.Range("D21").Value = 7
$ws.Range("E21").Value = 14.285714285714
$ws.Range("F21").Value = 43
$ws.Range("G21").Value = 32
$ws.Range("H21").Value = 34.375
$ws.Range("I21").Value = 425
$ws.Range("J21").Value = 390
$ws.Range("K21").Value = 8.974358974358
$ws.Range("L21").Value = 2.163461538461
$ws.Range("M21").Value = -8.798283261802
$ws.Range("N21").Value = -68.565088757396
$ws.Range("L22").Value = -66.666666666666
$ws.Range("F23").Value = 6
$ws.Range("H23").Value = 50
$ws.Range("I23").Value = 82
$ws.Range("K23").Value = -4.651162790697
$ws.Range("L23").Value = -10.869565217391
$ws.Range("M23").Value = 28.125
$ws.Range("D24").Value = 16
$ws.Range("E24").Value = -31.25
$ws.Range("F24").Value = 45
$ws.Range("G24").Value = 62
$ws.Range("H24").Value = -27.419354838709
$ws.Range("I24").Value = 406
$ws.Range("J24").Value = 474
$ws.Range("K24").Value = -14.345991561181
$ws.Range("L24").Value = 18.713450292397
$ws.Range("M24").Value = 18.367346938775
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = -63.636363636363
$ws.Range("F25").Value = 18
$ws.Range("G25").Value = 31
$ws.Range("H25").Value = -41.935483870967
$ws.Range("I25").Value = 203
$ws.Range("J25").Value = 244
$ws.Range("K25").Value = -16.803278688524
$ws.Range("L25").Value = 125.555555555556
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 100
$ws.Range("F26").Value = 23
$ws.Range("G26").Value = 16
$ws.Range("H26").Value = 43.75
$ws.Range("I26").Value = 126
$ws.Range("J26").Value = 125
$ws.Range("K26").Value = 0.8
$ws.Range("L26").Value = -11.888111888111
$ws.Range("M26").Value = -40.566037735849
$ws.Range("D27").Value = 1
$ws.Range("D27").NumberFormat = "#,##0"
$ws.Range("E27").Value = -100
$ws.Range("E27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = -33.333333333333
$ws.Range("J27").Value = 9
$ws.Range("K27").Value = -11.111111111111
$ws.Range("F28").Value = 1
$ws.Range("H28").Value = 0
$ws.Range("G29").Value = 1
$ws.Range("N29").Value = -91.111111111111
$ws.Range("G30").Value = 1
$ws.Range("N30").Value = -89.189189189189
